# Implemented getting kafka relations.
# Update the "classFields" sheet (field name / field type) ordering for the
# com.macro.mall.demo.dto.PmsBrandDto rows (2-9) and swap the
# com.macro.mall.demo.controller.DemoController LOGGER / demoService rows (11-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # classFields

# PmsBrandDto fields (rows 2-9): reorder field name / type pairs
$ws.Range("B2").Value = "sort"
$ws.Range("D2").Value = "java.lang.Integer"

$ws.Range("B3").Value = "showStatus"
$ws.Range("D3").Value = "java.lang.Integer"

$ws.Range("B4").Value = "name"
$ws.Range("D4").Value = "java.lang.String"

$ws.Range("B5").Value = "factoryStatus"
$ws.Range("D5").Value = "java.lang.Integer"

$ws.Range("B6").Value = "firstLetter"
$ws.Range("D6").Value = "java.lang.String"

$ws.Range("B7").Value = "logo"
$ws.Range("D7").Value = "java.lang.String"

$ws.Range("B8").Value = "brandStory"
$ws.Range("D8").Value = "java.lang.String"

$ws.Range("B9").Value = "bigPic"
$ws.Range("D9").Value = "java.lang.String"

# DemoController fields (rows 11-12): swap LOGGER / demoService
$ws.Range("B11").Value = "demoService"
$ws.Range("D11").Value = "com.macro.mall.demo.service.DemoService"

$ws.Range("B12").Value = "LOGGER"
$ws.Range("D12").Value = "org.slf4j.Logger"
